$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add P1 = 14, Q1 = 15, matching the existing header style ---
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25: swap I<->K and M<->O, then add P/Q = 2 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
